$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 8: Dokumentation
$ws.Range("A8").Value = 6
$ws.Range("B8").Value = "Dokumentation"
$ws.Range("F8").Value = 5
$ws.Rows.Item(8).RowHeight = 12.75

# New row 9: Inlärning
$ws.Range("A9").Value = 7
$ws.Range("B9").Value = "Inlärning"
$ws.Range("F9").Value = 6
$ws.Rows.Item(9).RowHeight = 12.75

# Widen column A slightly
$ws.Columns.Item(1).ColumnWidth = 3.6

# Update the selection to match the diff
$ws.Range("B20").Select()
